$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Helper: write a value to a cell while preserving the cell's "text"
# storage type (i.e. so a 4-decimal numeric-looking string such as
# "1.5627" is kept as text, matching the original workbook layout,
# instead of silently being converted into a numeric cell by Excel).
# -----------------------------------------------------------------
function Set-TextValue {
    param($cell, [string]$value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# ===================================================================
# Sheet: P_valores  (pairwise p-values matrix)
# ===================================================================
$wsP = $wb.Worksheets.Item("P_valores")

$wsP.Range("E2").Value = 0.005305285977035057
$wsP.Range("F2").Value = 0.004131548874835911
$wsP.Range("H2").Value = 0.00412467555736673
$wsP.Range("J2").Value = 0.005647142000966943

$wsP.Range("E3").Value = 0.02393143720546975
$wsP.Range("F3").Value = 0.01871609137959562
$wsP.Range("H3").Value = 0.008875912317723778
$wsP.Range("J3").Value = 0.008233965610212612

$wsP.Range("E4").Value = 0.9485800820697468
$wsP.Range("F4").Value = 0.9103669270129748
$wsP.Range("H4").Value = 0.03306027981392456
$wsP.Range("J4").Value = 0.01347114938333038

$wsP.Range("B5").Value = 0.005305285977035057
$wsP.Range("C5").Value = 0.02393143720546975
$wsP.Range("D5").Value = 0.9485800820697468
$wsP.Range("F5").Value = 0.4863314953090097
$wsP.Range("G5").Value = 0.3840565803825529
$wsP.Range("H5").Value = 0.006038550104455975
$wsP.Range("I5").Value = 0.005629114665586243
$wsP.Range("J5").Value = 0.005904426321486245

$wsP.Range("B6").Value = 0.004131548874835911
$wsP.Range("C6").Value = 0.01871609137959562
$wsP.Range("D6").Value = 0.9103669270129748
$wsP.Range("E6").Value = 0.4863314953090097
$wsP.Range("G6").Value = 0.4488142296517927
$wsP.Range("H6").Value = 0.007730055059711249
$wsP.Range("I6").Value = 0.006341787075952121
$wsP.Range("J6").Value = 0.006524913545280064

$wsP.Range("E7").Value = 0.3840565803825529
$wsP.Range("F7").Value = 0.4488142296517927
$wsP.Range("H7").Value = 0.02546020839792718
$wsP.Range("J7").Value = 0.0009665947134080977

$wsP.Range("B8").Value = 0.00412467555736673
$wsP.Range("C8").Value = 0.008875912317723778
$wsP.Range("D8").Value = 0.03306027981392456
$wsP.Range("E8").Value = 0.006038550104455975
$wsP.Range("F8").Value = 0.007730055059711249
$wsP.Range("G8").Value = 0.02546020839792718
$wsP.Range("I8").Value = 0.01493134748401359
$wsP.Range("J8").Value = 0.01086650606326289

$wsP.Range("E9").Value = 0.005629114665586243
$wsP.Range("F9").Value = 0.006341787075952121
$wsP.Range("H9").Value = 0.01493134748401359
$wsP.Range("J9").Value = 0.01376997561758664

$wsP.Range("B10").Value = 0.005647142000966943
$wsP.Range("C10").Value = 0.008233965610212612
$wsP.Range("D10").Value = 0.01347114938333038
$wsP.Range("E10").Value = 0.005904426321486245
$wsP.Range("F10").Value = 0.006524913545280064
$wsP.Range("G10").Value = 0.0009665947134080977
$wsP.Range("H10").Value = 0.01086650606326289
$wsP.Range("I10").Value = 0.01376997561758664

# ===================================================================
# Sheet: Estadisticos_HLN_DM  (pairwise HLN-DM statistics matrix)
# ===================================================================
$wsH = $wb.Worksheets.Item("Estadisticos_HLN_DM")

$wsH.Range("E2").Value = -2.979011395689851
$wsH.Range("F2").Value = -3.075301860208742
$wsH.Range("H2").Value = -3.075938670891716
$wsH.Range("J2").Value = -2.954753774119182

$wsH.Range("E3").Value = -2.364024560027169
$wsH.Range("F3").Value = -2.469350628913434
$wsH.Range("H3").Value = -2.776326266125144
$wsH.Range("J3").Value = -2.806305432527586

$wsH.Range("E4").Value = 0.06496786511420824
$wsH.Range("F4").Value = -0.1134171223842102
$wsH.Range("H4").Value = -2.221754490375471
$wsH.Range("J4").Value = -2.606822986779918

$wsH.Range("B5").Value = 2.979011395689851
$wsH.Range("C5").Value = 2.364024560027169
$wsH.Range("D5").Value = -0.06496786511420824
$wsH.Range("F5").Value = -0.7038296828992813
$wsH.Range("G5").Value = -0.8818376846265282
$wsH.Range("H5").Value = -2.928623257780381
$wsH.Range("I5").Value = -2.955997956507879
$wsH.Range("J5").Value = -2.937393079054316

$wsH.Range("B6").Value = 3.075301860208742
$wsH.Range("C6").Value = 2.469350628913434
$wsH.Range("D6").Value = 0.1134171223842102
$wsH.Range("E6").Value = 0.7038296828992813
$wsH.Range("G6").Value = -0.7662515165628037
$wsH.Range("H6").Value = -2.831409553940748
$wsH.Range("I6").Value = -2.90945244358324
$wsH.Range("J6").Value = -2.898288175503464

$wsH.Range("E7").Value = 0.8818376846265282
$wsH.Range("F7").Value = 0.7662515165628037
$wsH.Range("H7").Value = -2.337111901495963
$wsH.Range("J7").Value = -3.612944267251107

$wsH.Range("B8").Value = 3.075938670891716
$wsH.Range("C8").Value = 2.776326266125144
$wsH.Range("D8").Value = 2.221754490375471
$wsH.Range("E8").Value = 2.928623257780381
$wsH.Range("F8").Value = 2.831409553940748
$wsH.Range("G8").Value = 2.337111901495963
$wsH.Range("I8").Value = -2.564193228957207
$wsH.Range("J8").Value = -2.694754620437258

$wsH.Range("E9").Value = 2.955997956507879
$wsH.Range("F9").Value = 2.90945244358324
$wsH.Range("H9").Value = 2.564193228957207
$wsH.Range("J9").Value = -2.597763446613827

$wsH.Range("B10").Value = 2.954753774119182
$wsH.Range("C10").Value = 2.806305432527586
$wsH.Range("D10").Value = 2.606822986779918
$wsH.Range("E10").Value = 2.937393079054316
$wsH.Range("F10").Value = 2.898288175503464
$wsH.Range("G10").Value = 3.612944267251107
$wsH.Range("H10").Value = 2.694754620437258
$wsH.Range("I10").Value = 2.597763446613827

# ===================================================================
# Sheet: Resumen_Modelos  (summary table - ECRPS stats stored as text)
# ===================================================================
$wsR = $wb.Worksheets.Item("Resumen_Modelos")

Set-TextValue $wsR.Range("F5")  "1.5627"
Set-TextValue $wsR.Range("F6")  "1.5887"
Set-TextValue $wsR.Range("F8")  "2.4549"
Set-TextValue $wsR.Range("F10") "4.2755"

Set-TextValue $wsR.Range("G5")  "3.5113"
Set-TextValue $wsR.Range("G6")  "3.3264"
Set-TextValue $wsR.Range("G8")  "5.4040"
Set-TextValue $wsR.Range("G10") "12.2657"

Set-TextValue $wsR.Range("H5")  "2.2469"
Set-TextValue $wsR.Range("H6")  "2.0937"
Set-TextValue $wsR.Range("H8")  "2.2013"
Set-TextValue $wsR.Range("H10") "2.8688"

$wb.Save()
